# Update column F ("F" = 报名/浏览量 count column) values on each worksheet
# as per the commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 251
$ws.Cells.Item(4, 6).Value = 847
$ws.Cells.Item(6, 6).Value = 418
$ws.Cells.Item(7, 6).Value = 606
$ws.Cells.Item(8, 6).Value = 225
$ws.Cells.Item(10, 6).Value = 360
$ws.Cells.Item(11, 6).Value = 156
$ws.Cells.Item(12, 6).Value = 698
$ws.Cells.Item(13, 6).Value = 93
$ws.Cells.Item(14, 6).Value = 1828
$ws.Cells.Item(15, 6).Value = 366
$ws.Cells.Item(16, 6).Value = 3747
$ws.Cells.Item(17, 6).Value = 373
$ws.Cells.Item(18, 6).Value = 497
$ws.Cells.Item(19, 6).Value = 7
$ws.Cells.Item(20, 6).Value = 59
$ws.Cells.Item(21, 6).Value = 147

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 225
$ws.Cells.Item(7, 6).Value = 480
$ws.Cells.Item(13, 6).Value = 95
$ws.Cells.Item(21, 6).Value = 23

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 5365
$ws.Cells.Item(4, 6).Value = 282

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 5365
$ws.Cells.Item(6, 6).Value = 282
$ws.Cells.Item(7, 6).Value = 251
$ws.Cells.Item(8, 6).Value = 225
$ws.Cells.Item(12, 6).Value = 480
$ws.Cells.Item(13, 6).Value = 480
$ws.Cells.Item(14, 6).Value = 847
$ws.Cells.Item(18, 6).Value = 418
$ws.Cells.Item(19, 6).Value = 606
$ws.Cells.Item(20, 6).Value = 225
$ws.Cells.Item(23, 6).Value = 360
$ws.Cells.Item(24, 6).Value = 156
$ws.Cells.Item(27, 6).Value = 698
$ws.Cells.Item(28, 6).Value = 93
$ws.Cells.Item(29, 6).Value = 95
$ws.Cells.Item(30, 6).Value = 1828
$ws.Cells.Item(31, 6).Value = 366
$ws.Cells.Item(32, 6).Value = 3748
$ws.Cells.Item(34, 6).Value = 373
$ws.Cells.Item(35, 6).Value = 497
$ws.Cells.Item(36, 6).Value = 7
$ws.Cells.Item(37, 6).Value = 59
$ws.Cells.Item(39, 6).Value = 147
$ws.Cells.Item(46, 6).Value = 23
